# This script applies a weekly data update to the "Pera" (pear) price sheet.
# A new week of data (3 rows) is inserted right before the existing block for
# Fecha=44448 (rows 276-278 in the original sheet), pushing the rest of the
# data for this product down by 3 rows. The new rows report prices observed
# for Fecha=45041 (the most recent week) in the Región de O'Higgins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 276:278. Excel copies formatting (including the
# date number format on column D) from the row above, which is exactly what
# we need for the new rows.
$ws.Rows("276:278").Insert()

# Shared values for the 3 new rows (same market/product/variety, new date).
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$fecha     = 45041
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100104
$producto  = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$variedad  = "Packham's Triumph"

# Row 276: Especial
$ws.Range("A276").Value = 7
$ws.Range("B276").Value = $mercado
$ws.Range("C276").Value = $region
$ws.Range("D276").Value = $fecha
$ws.Range("E276").Value = $codreg
$ws.Range("F276").Value = $tipo
$ws.Range("G276").Value = $productoId
$ws.Range("H276").Value = $producto
$ws.Range("I276").Value = $categoriaId
$ws.Range("J276").Value = $categoria
$ws.Range("K276").Value = $variedad
$ws.Range("L276").Value = "Especial"
$ws.Range("M276").Value = 80
$ws.Range("N276").Value = 12000
$ws.Range("O276").Value = 12000
$ws.Range("P276").Value = 12000
$ws.Range("Q276").Value = "$/bandeja 18 kilos granel"
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value = 667
$ws.Range("T276").Value = 18

# Row 277: Primera
$ws.Range("A277").Value = 7
$ws.Range("B277").Value = $mercado
$ws.Range("C277").Value = $region
$ws.Range("D277").Value = $fecha
$ws.Range("E277").Value = $codreg
$ws.Range("F277").Value = $tipo
$ws.Range("G277").Value = $productoId
$ws.Range("H277").Value = $producto
$ws.Range("I277").Value = $categoriaId
$ws.Range("J277").Value = $categoria
$ws.Range("K277").Value = $variedad
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 80
$ws.Range("N277").Value = 10000
$ws.Range("O277").Value = 10000
$ws.Range("P277").Value = 10000
$ws.Range("Q277").Value = "$/bandeja 18 kilos granel"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 556
$ws.Range("T277").Value = 18

# Row 278: Segunda
$ws.Range("A278").Value = 7
$ws.Range("B278").Value = $mercado
$ws.Range("C278").Value = $region
$ws.Range("D278").Value = $fecha
$ws.Range("E278").Value = $codreg
$ws.Range("F278").Value = $tipo
$ws.Range("G278").Value = $productoId
$ws.Range("H278").Value = $producto
$ws.Range("I278").Value = $categoriaId
$ws.Range("J278").Value = $categoria
$ws.Range("K278").Value = $variedad
$ws.Range("L278").Value = "Segunda"
$ws.Range("M278").Value = 50
$ws.Range("N278").Value = 9000
$ws.Range("O278").Value = 9000
$ws.Range("P278").Value = 9000
$ws.Range("Q278").Value = "$/bandeja 18 kilos granel"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 500
$ws.Range("T278").Value = 18
